$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of test data (Tanggal, No. Docket, Kontraktor, Nama Proyek, Teknisi,
# Nama Mutu, Slump, TM, Waktu Kirim) appended below the existing 19 data rows.
$newRows = @(
    @("2025-12-04", "DOCKET/PTB/12-2025/26867", "HK - JAKON JO / DECKSLAB A1 - P OP 5 / CLASS B-1 / AD PTB", "Pembangunan Jalan Tol  Akses Patimban Paket 3", "NANDA", "Class B-1 FA", " Slump 12.0 +2.0/-2.0", "TM 496", "2025-12-04 08:55:21"),
    @("2025-12-04", "DOCKET/PTB/12-2025/26869", "HK - JAKON JO / DECKSLAB A1 - P OP 5 / CLASS B-1 / AD PTB", "Pembangunan Jalan Tol  Akses Patimban Paket 3", "NANDA", "Class B-1 FA", " Slump 12.0 +2.0/-2.0", "TM 601", "2025-12-04 09:36:50"),
    @("2025-12-04", "DOCKET/PTB/12-2025/26883", "HK - JAKON JO / DECKSLAB A1 - P OP 5 / CLASS B-1 / AD PTB", "Pembangunan Jalan Tol  Akses Patimban Paket 3", "AGUS", "Class B-1 FA", " Slump 12.0 +2.0/-2.0", "TM 965", "2025-12-04 19:57:37"),
    @("2025-12-04", "DOCKET/PTB/12-2025/26864", "HK - JAKON JO / LC A2 OP 7 / CLASS E / AD PTB", "Pembangunan Jalan Tol  Akses Patimban Paket 3", "AGUS", "Class E-1 FA", " Slump 12.0 +2.0/-2.0", "TM 690", "2025-12-04 08:25:26"),
    @("2025-12-04", "DOCKET/PTB/12-2025/26865", "HK - JAKON JO / LC A2 OP 7 / CLASS E / AD PTB", "Pembangunan Jalan Tol  Akses Patimban Paket 3", "AGUS", "Class E-1 FA", " Slump 12.0 +2.0/-2.0", "TM 944", "2025-12-04 08:31:24"),
    @("2025-12-04", "DOCKET/PTB/12-2025/26880", "PT SRA / K 350 / AD PTB", "Proyek Pekerjaan Fasilitas Workshop PGT Patimban", "SUGENG", "K-350 NFA", " Slump 12.0 +2.0/-2.0", "TM 813", "2025-12-04 14:01:51"),
    @("2025-12-04", "DOCKET/PTB/12-2025/26875", "WASKITA - ABP JO / BARRIER GRID AE2 - PE1' STA 25+949 - 25+896 SISI MEDIAN & LUAR / CLASS B-1 / AD PTB", "Patimban Access Toll Road P02 (JOI 60%)", "AGUS", "Class B-1 NFA", " Slump 12.0 +2.0/-2.0", "TM 965", "2025-12-04 11:41:00"),
    @("2025-12-04", "DOCKET/PTB/12-2025/26881", "WASKITA - ABP JO / BARRIER GRID AE2 - PE1' STA 25+949 - 25+896 SISI MEDIAN & LUAR / CLASS B-1 / AD PTB", "Patimban Access Toll Road P02 (JOI 60%)", "AGUS", "Class B-1 NFA", " Slump 12.0 +2.0/-2.0", "TM 827", "2025-12-04 14:11:17")
)

$startRow = 20
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 1; $c -le 9; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = $row[$c - 1]
        $cell.Style = $ws.Cells.Item($r - 1, $c).Style
    }
}

$lastRow = $startRow + $newRows.Count - 1

# Update auto filter / selection to reflect new data extent.
$ws.Range("A1:I1").AutoFilter() | Out-Null

$selRow = $lastRow + 1
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A$selRow`:XFD$selRow").Select() | Out-Null
